$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 (S.NO = 9): Chroma + LangChain + RAG + Text
$ws.Range("C10").Value = "Chroma + LangChain + RAG + Text"
$ws.Range("D10").Value = "LangChain + Gemini"
$ws.Range("E10").Value = "Done"

# Row 11 (S.NO = 10): Chroma + LangChain + RAG + PDF
$ws.Range("C11").Value = "Chroma + LangChain + RAG + PDF"
$ws.Range("D11").Value = "LangChain + Gemini"
$ws.Range("E11").Value = "Done"

# Apply the "Done" style (same as other Status cells, e.g. E9) to the new Status cells
$ws.Range("E9").Copy() | Out-Null
$ws.Range("E10:E11").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Update the active selection to G11 (as reflected in the saved sheet view)
$ws.Range("G11").Select() | Out-Null
